$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.036.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.73%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.885.56'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.07%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7363'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.53%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.02%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9993'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.12%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3159'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.42%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07163'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.13%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.72'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.74%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08322'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.04%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7557'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.25%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.898.99'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.78%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.392'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.69%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.60'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.04%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.141'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.41%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.043.50'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.80%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '248.91'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.75%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.55'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.88%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007849'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.47%  '

# Row 21
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9987'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.03%  '

# Row 22
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.131.31'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.10%  '

# Row 23
$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.889'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.25%  '

# Row 24
$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.10%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1565'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.89%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.272'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.53%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.70'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.47%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.69'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.14%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.045'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.03%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.475'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.01%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.570'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.83%  '

# Row 32
$ws.Range("E32").Value = '  -0.21%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.191'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.80%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05322'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.35%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.248'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.97%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7691'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.30%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9967'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.97%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.721'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.59%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01957'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.82%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.758'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.31%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4581'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.86%  '

# Row 42
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.030'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.95%  '

# Row 43
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.087.48'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.84%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8778'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.07%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '72.31'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.14%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '104.29'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.69%  '

# Row 47
$ws.Range("E47").Value = '  +0.00%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.854'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.22%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.562'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.95%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.548'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.44%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.037.68'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.10%  '
